# Atualizacao 16 nov 2020
# Appends the newest daily/monthly readings to the "Rio Negro" station
# workbook:
#   - "Mensal" sheet: one new monthly sample (row 14)
#   - "Diario" sheet: fifteen new daily samples (rows 368-382, 2020-11-01..15)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Mensal" -> add row 14
# ---------------------------------------------------------------------
$wsMensal = $wb.Worksheets.Item("Mensal")

# Carry the date-column formatting (bold, centered, bordered, yyyy-mm-dd)
# down from the last existing row so the new cell picks up the same style
# index the rest of column A uses.
$wsMensal.Range("A13").Copy() | Out-Null
$wsMensal.Range("A14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsMensal.Range("A14").Value = 44150
$wsMensal.Range("B14").Value = 53.71
$wsMensal.Range("C14").Value = 149.37
$wsMensal.Range("D14").Value = -64.04000000000001

# ---------------------------------------------------------------------
# Sheet "Diario" -> add rows 368-382
# ---------------------------------------------------------------------
$wsDiario = $wb.Worksheets.Item("Diario")

# Same trick for the whole new block of date cells at once.
$wsDiario.Range("A367").Copy() | Out-Null
$wsDiario.Range("A368:A382").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$diarioData = @(
    @(368, 44136, 48.6,               149.37, -67.45999999999999),
    @(369, 44137, 50.76,              149.37, -66.02),
    @(370, 44138, 53.54,              149.37, -64.15000000000001),
    @(371, 44139, 51.76,              149.37, -65.34999999999999),
    @(372, 44140, 47.41,              149.37, -68.26000000000001),
    @(373, 44141, 44.78,              149.37, -70.02),
    @(374, 44142, 40.89,              149.37, -72.63),
    @(375, 44143, 38.01,              149.37, -74.55),
    @(376, 44144, 35.18,              149.37, -76.45),
    @(377, 44145, 36.29,              149.37, -75.70999999999999),
    @(378, 44146, 41.03,              149.37, -72.53),
    @(379, 44147, 54.06,              149.37, -63.81),
    @(380, 44148, 96.73999999999999,  149.37, -35.23),
    @(381, 44149, 88.93000000000001,  149.37, -40.46),
    @(382, 44150, 77.66,              149.37, -48.01)
)

foreach ($row in $diarioData) {
    $r = $row[0]
    $wsDiario.Cells.Item($r, 1).Value = $row[1]
    $wsDiario.Cells.Item($r, 2).Value = $row[2]
    $wsDiario.Cells.Item($r, 3).Value = $row[3]
    $wsDiario.Cells.Item($r, 4).Value = $row[4]
}
